$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.531.07"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.511.84"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "2.903.31"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "2.494.70"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "48.346.23"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").Value = "0.0₃0935"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.37%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.143"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0771"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0301"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "1.991.09"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("E47").Value = "  +5.61%  "
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
